$d = $word.ActiveDocument

# --- Change 1: " (коефіцієнт подібності)" -> " (косинус подібності)" ---
# Locate the substring "оефіцієнт" (everything after the leading "к") and
# replace it with "осинус", forcing the run to split into three runs that
# share identical formatting (matches the target OOXML structure).
$r1 = $d.Content
$r1.Find.Execute("оефіцієнт", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start1 = $r1.Start
$r1.Text = "осинус"
$mid1 = $d.Range($start1, $start1 + 6)
$mid1.Bold = 1
$mid1.Bold = 0

# --- Change 2: "(з коефіцієнтом подібності більше 0,8) ..." -> "(з косинусом подібності більше 0,8) ..." ---
# Locate "коефіцієнтом " (with trailing space) and replace it with
# "косинусом " the same way.
$r2 = $d.Content
$r2.Find.Execute("коефіцієнтом ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start2 = $r2.Start
$r2.Text = "косинусом "
$mid2 = $d.Range($start2, $start2 + 10)
$mid2.Bold = 1
$mid2.Bold = 0
